$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the USD Amount value in T2 (284999 -> 384005)
$ws.Range("T2").Value = 384005

# Update the sheet's active cell / selection to R11
$ws.Range("R11").Select()
